$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update headline metrics after trade #100 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.08   # Current Capital
$summary.Range("B4").Value = 0.09      # Total P&L $
$summary.Range("B6").Value = 100       # Total Trades
$summary.Range("B8").Value = 41        # Losing Trades
$summary.Range("B9").Value = 41        # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.08     # Capital
$status.Range("D4").Value = 100        # Trades
$status.Range("E4").Value = 0.09       # P&L $
$status.Range("F4").Value = 0.08       # P&L %
$status.Range("G4").Value = 41         # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade #100 (row 101) to "All Trades" and
# "MarketMaking" sheets - both sheets mirror the same trade log.
# ---------------------------------------------------------------------------
$newRow = @{
    A = 100
    B = "2026-02-17"
    C = "09:17:47"
    D = "MarketMaking"
    E = "UP"
    F = 0.14
    G = 0.11
    H = "CLOSED"
    I = -21.4286
    J = -0.03
    K = 100.08
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A101").Value = $newRow.A

    # Keep the date / time strings as plain text instead of letting Excel
    # auto-convert them into date/time serial values.
    $ws.Range("B101").NumberFormat = "@"
    $ws.Range("B101").Value = $newRow.B
    $ws.Range("C101").NumberFormat = "@"
    $ws.Range("C101").Value = $newRow.C

    $ws.Range("D101").Value = $newRow.D
    $ws.Range("E101").Value = $newRow.E
    $ws.Range("F101").Value = $newRow.F
    $ws.Range("G101").Value = $newRow.G
    $ws.Range("H101").Value = $newRow.H
    $ws.Range("I101").Value = $newRow.I
    $ws.Range("J101").Value = $newRow.J
    $ws.Range("K101").Value = $newRow.K
    $ws.Range("L101").Value = $newRow.L
    $ws.Range("M101").Value = $newRow.M
    $ws.Range("N101").Value = $newRow.N
    $ws.Range("O101").Value = $newRow.O
    $ws.Range("P101").Value = $newRow.P
    $ws.Range("Q101").Value = $newRow.Q
}
